$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 40.83537666666667
$ws.Range("H2").Value = 122.50613
$ws.Range("I2").Value = 0.9274830900091532
$ws.Range("J2").Value = 0.9274830900091531
$ws.Range("M2").Value = 127.3992563333333
$ws.Range("N2").Value = 382.197769
$ws.Range("O2").Value = 0.4838549810199306
$ws.Range("P2").Value = 0.4838549810199307
$ws.Range("Q2").Value = 5202.396619424885
$ws.Range("R2").Value = 46821.56957482397
$ws.Range("S2").Value = 0.4487673129126854
$ws.Range("T2").Value = 0.4487673129126855

# Row 3
$ws.Range("G3").Value = 40.83537666666667
$ws.Range("H3").Value = 122.50613
$ws.Range("I3").Value = 0.9274830900091532
$ws.Range("J3").Value = 0.9274830900091531
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("N3").Value = 178.097596
$ws.Range("O3").Value = 0.2254681108101269
$ws.Range("P3").Value = 0.2254681108101269
$ws.Range("Q3").Value = 2424.227472029275
$ws.Range("R3").Value = 21818.04724826348
$ws.Range("S3").Value = 0.2091178601127026
$ws.Range("T3").Value = 0.2091178601127026

# Row 4
$ws.Range("G4").Value = 40.83537666666667
$ws.Range("H4").Value = 122.50613
$ws.Range("I4").Value = 0.9274830900091532
$ws.Range("J4").Value = 0.9274830900091531
$ws.Range("M4").Value = 16.63275166666667
$ws.Range("N4").Value = 49.898255
$ws.Range("O4").Value = 0.06317022542837675
$ws.Range("P4").Value = 0.06317022542837675
$ws.Range("Q4").Value = 679.2046793114612
$ws.Range("R4").Value = 6112.842113803151
$ws.Range("S4").Value = 0.05858931587688565
$ws.Range("T4").Value = 0.05858931587688564

# Row 5
$ws.Range("G5").Value = 40.83537666666667
$ws.Range("H5").Value = 122.50613
$ws.Range("I5").Value = 0.9274830900091532
$ws.Range("J5").Value = 0.9274830900091531
$ws.Range("M5").Value = 59.90262233333334
$ws.Range("N5").Value = 179.707867
$ws.Range("O5").Value = 0.2275066827415657
$ws.Range("P5").Value = 0.2275066827415658
$ws.Range("Q5").Value = 2446.146146302746
$ws.Range("R5").Value = 22015.31531672471
$ws.Range("S5").Value = 0.2110086011068795
$ws.Range("T5").Value = 0.2110086011068795

# Row 6
$ws.Range("I6").Value = 0.03813623414934058
$ws.Range("J6").Value = 0.03813623414934057
$ws.Range("M6").Value = 127.3992563333333
$ws.Range("N6").Value = 382.197769
$ws.Range("O6").Value = 0.4838549810199306
$ws.Range("P6").Value = 0.4838549810199307
$ws.Range("Q6").Value = 213.9120569995161
$ws.Range("R6").Value = 1925.208512995645
$ws.Range("S6").Value = 0.01845240685050082
$ws.Range("T6").Value = 0.01845240685050082

# Row 7
$ws.Range("I7").Value = 0.03813623414934058
$ws.Range("J7").Value = 0.03813623414934057
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("N7").Value = 178.097596
$ws.Range("O7").Value = 0.2254681108101269
$ws.Range("P7").Value = 0.2254681108101269
$ws.Range("Q7").Value = 99.6793445621311
$ws.Range("R7").Value = 897.1141010591799
$ws.Range("S7").Value = 0.008598504667064466
$ws.Range("T7").Value = 0.008598504667064466

# Row 8
$ws.Range("I8").Value = 0.03813623414934058
$ws.Range("J8").Value = 0.03813623414934057
$ws.Range("M8").Value = 16.63275166666667
$ws.Range("N8").Value = 49.898255
$ws.Range("O8").Value = 0.06317022542837675
$ws.Range("P8").Value = 0.06317022542837675
$ws.Range("Q8").Value = 27.92752661969723
$ws.Range("R8").Value = 251.347739577275
$ws.Range("S8").Value = 0.002409074508203204
$ws.Range("T8").Value = 0.002409074508203204

# Row 9
$ws.Range("I9").Value = 0.03813623414934058
$ws.Range("J9").Value = 0.03813623414934057
$ws.Range("M9").Value = 59.90262233333334
$ws.Range("N9").Value = 179.707867
$ws.Range("O9").Value = 0.2275066827415657
$ws.Range("P9").Value = 0.2275066827415658
$ws.Range("Q9").Value = 100.5805962435261
$ws.Range("R9").Value = 905.2253661917351
$ws.Range("S9").Value = 0.008676248123572093
$ws.Range("T9").Value = 0.008676248123572093

# Row 10
$ws.Range("G10").Value = 1.503819
$ws.Range("H10").Value = 4.511457
$ws.Range("I10").Value = 0.03415584247746153
$ws.Range("J10").Value = 0.03415584247746152
$ws.Range("M10").Value = 127.3992563333333
$ws.Range("N10").Value = 382.197769
$ws.Range("O10").Value = 0.4838549810199306
$ws.Range("P10").Value = 0.4838549810199307
$ws.Range("Q10").Value = 191.585422259937
$ws.Range("R10").Value = 1724.268800339433
$ws.Range("S10").Value = 0.01652647451365189
$ws.Range("T10").Value = 0.01652647451365189

# Row 11
$ws.Range("G11").Value = 1.503819
$ws.Range("H11").Value = 4.511457
$ws.Range("I11").Value = 0.03415584247746153
$ws.Range("J11").Value = 0.03415584247746152
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("N11").Value = 178.097596
$ws.Range("O11").Value = 0.2254681108101269
$ws.Range("P11").Value = 0.2254681108101269
$ws.Range("Q11").Value = 89.27551623970798
$ws.Range("R11").Value = 803.4796461573719
$ws.Range("S11").Value = 0.007701053276521533
$ws.Range("T11").Value = 0.007701053276521533

# Row 12
$ws.Range("G12").Value = 1.503819
$ws.Range("H12").Value = 4.511457
$ws.Range("I12").Value = 0.03415584247746153
$ws.Range("J12").Value = 0.03415584247746152
$ws.Range("M12").Value = 16.63275166666667
$ws.Range("N12").Value = 49.898255
$ws.Range("O12").Value = 0.06317022542837675
$ws.Range("P12").Value = 0.06317022542837675
$ws.Range("Q12").Value = 25.012647978615
$ws.Range("R12").Value = 225.113831807535
$ws.Range("S12").Value = 0.002157632268997371
$ws.Range("T12").Value = 0.00215763226899737

# Row 13
$ws.Range("G13").Value = 1.503819
$ws.Range("H13").Value = 4.511457
$ws.Range("I13").Value = 0.03415584247746153
$ws.Range("J13").Value = 0.03415584247746152
$ws.Range("M13").Value = 59.90262233333334
$ws.Range("N13").Value = 179.707867
$ws.Range("O13").Value = 0.2275066827415657
$ws.Range("P13").Value = 0.2275066827415658
$ws.Range("Q13").Value = 90.08270161469102
$ws.Range("R13").Value = 810.744314532219
$ws.Range("S13").Value = 0.007770682418290735
$ws.Range("T13").Value = 0.007770682418290735

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.009899
$ws.Range("H14").Value = 0.029697
$ws.Range("I14").Value = 0.0002248333640447365
$ws.Range("J14").Value = 0.0002248333640447365
$ws.Range("M14").Value = 127.3992563333333
$ws.Range("N14").Value = 382.197769
$ws.Range("O14").Value = 0.4838549810199306
$ws.Range("P14").Value = 0.4838549810199307
$ws.Range("Q14").Value = 1.261125238443667
$ws.Range("R14").Value = 11.350127145993
$ws.Range("S14").Value = 0.0001087867430925132
$ws.Range("T14").Value = 0.0001087867430925132

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.009899
$ws.Range("H15").Value = 0.029697
$ws.Range("I15").Value = 0.0002248333640447365
$ws.Range("J15").Value = 0.0002248333640447365
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("N15").Value = 178.097596
$ws.Range("O15").Value = 0.2254681108101269
$ws.Range("P15").Value = 0.2254681108101269
$ws.Range("Q15").Value = 0.5876627009346665
$ws.Range("R15").Value = 5.288964308412
$ws.Range("S15").Value = 0.00005069275383825225
$ws.Range("T15").Value = 0.00005069275383825225

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.009899
$ws.Range("H16").Value = 0.029697
$ws.Range("I16").Value = 0.0002248333640447365
$ws.Range("J16").Value = 0.0002248333640447365
$ws.Range("M16").Value = 16.63275166666667
$ws.Range("N16").Value = 49.898255
$ws.Range("O16").Value = 0.06317022542837675
$ws.Range("P16").Value = 0.06317022542837675
$ws.Range("Q16").Value = 0.1646476087483333
$ws.Range("R16").Value = 1.481828478735
$ws.Range("S16").Value = 0.0000142027742905263
$ws.Range("T16").Value = 0.0000142027742905263

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.009899
$ws.Range("H17").Value = 0.029697
$ws.Range("I17").Value = 0.0002248333640447365
$ws.Range("J17").Value = 0.0002248333640447365
$ws.Range("M17").Value = 59.90262233333334
$ws.Range("N17").Value = 179.707867
$ws.Range("O17").Value = 0.2275066827415657
$ws.Range("P17").Value = 0.2275066827415658
$ws.Range("Q17").Value = 0.5929760584776668
$ws.Range("R17").Value = 5.336784526299001
$ws.Range("S17").Value = 0.00005115109282344483
$ws.Range("T17").Value = 0.00005115109282344484
